# "Updated the diagram and module status"
#
# - Removes the "state.preloader" module row.
# - Adds a new "app.fancy-slider.assets-downloader" module row (right after
#   app.fancy-slider.animations, before app.fancy-slider.blur).
# - Expands the note for common.gsap-lite to mention TimelineLite / Power eases.
# - Widens column A and resets the sheet view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "app.fancy-slider.assets-downloader" right before the
# current row 4 (app.fancy-slider.blur). Excel copies formatting from the row
# above automatically, matching the "Good" style already used by B/C/D here.
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value2 = "app.fancy-slider.assets-downloader"
$ws.Range("E4").Value2 = "Foloseste preloader'ul pentru a crea o lista de descarcare. Cand preloaderul termina descarcarea, acesta mai infasoara odata resursele intr-un loader PIXI pentru a ne asigura ca acestea sunt si in cache'ul pixi."

# The old "state.preloader" row has shifted from row 19 down to row 20 because
# of the insertion above; delete it entirely.
$ws.Rows.Item(20).Delete()

# Update the common.gsap-lite note (now on row 13 after the insert/delete).
$ws.Range("E13").Value2 = "Incapsuleaza TweenLite, TimelineLite si Easing'urile (Power1, Power2, Power3, Power4)"

# Widen column A to display width 36 (ColumnWidth uses character units that
# are offset by 5/6 of a character from the stored sheet width).
$ws.Columns.Item(1).ColumnWidth = 35.16666666666666

# Reset the view: 100% zoom and selection on G8 (matching the saved view).
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("G8").Select()
